# Insert a new data row at row 42 (pushing existing rows 42-120 down to 43-121)
# and populate it with a new weekly price record for "Poroto granado".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(42).Insert()

$ws.Range("A42").Value = 8
$ws.Range("B42").Value = "Terminal La Palmera de La Serena"
$ws.Range("C42").Value = "Coquimbo"
$ws.Range("D42").Value = 44979
$ws.Range("E42").Value = 4
$ws.Range("F42").Value = 100112030
$ws.Range("G42").Value = "Poroto granado"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 400
$ws.Range("K42").Value = 37000
$ws.Range("L42").Value = 38000
$ws.Range("M42").Value = 37500
$ws.Range("N42").Value = "$/malla 25 kilos"
$ws.Range("O42").Value = "Provincia del Elquí"
$ws.Range("P42").Value = 1500
$ws.Range("Q42").Value = 25
$ws.Range("R42").Value = "Hortaliza"
